# Applies the Fruitilicious review edits:
#  - Title / headline text swapped to a shorter phrasing (appears twice).
#  - "What we like" bullet list: all four items replaced with new wording
#    (values shift to different bullets as part of the rewrite).
#  - "What we don't like" bullet list: both items replaced with new wording.
#  - Closing meta title + meta description paragraphs updated to match.
#
# Each replacement is scoped to its own paragraph's Range so that the
# substring relationships between some of the old/new phrases (e.g. the new
# text for one bullet containing the old text of another) cannot cause a
# document-wide Find/Replace to cascade into the wrong paragraph.

$d = $word.ActiveDocument

function Replace-InParagraph($index, $oldText, $newText) {
    $p = $d.Paragraphs.Item($index)
    $r = $p.Range
    $r.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
}

# Main page title (Heading1, paragraph 1)
Replace-InParagraph 1 "Play Fruitilicious Free: Review of Classic Online Slot Game" "Play Fruitilicious Slot Game for Free"

# "What we like" bullets (paragraphs 42-45)
Replace-InParagraph 42 "Wide range of betting options" "Straightforward gameplay with no special features or bonuses"
Replace-InParagraph 43 "Attractive fruit symbols and design" "Wide range of betting options for players of all levels"
Replace-InParagraph 44 "Potential for high payouts" "Potentially lucrative prizes"
Replace-InParagraph 45 "Simple and relaxing sound effects" "Fresh and attractive game design"

# "What we don't like" bullets (paragraphs 47-48)
Replace-InParagraph 47 "No special features or bonuses" "Lower than average RTP of 95%"
Replace-InParagraph 48 "Lower than average RTP" "No Wild or Scatter symbols"

# Bold "title" run near the end (paragraph 49)
Replace-InParagraph 49 "Play Fruitilicious Free: Review of Classic Online Slot Game" "Play Fruitilicious Slot Game for Free"

# Italic meta-description run (paragraph 50)
Replace-InParagraph 50 "Read our review of Fruitilicious, a classic online slot game with potential for high payouts. Play for free with a wide range of betting options." "Read our review of Fruitilicious, a straightforward slot game with potentially lucrative prizes. Play for free."
